$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A2 currently holds "76442711" as text; the edit turns it into a real number.
$ws.Range("A2").Value = 76442711

# New row 3: payment 76442781 (Cash) 2025-08-20T08:53:29
# A3 must stay textual (like the original A2 before this edit), even though it
# looks numeric, so force the cell to Text before writing it, then drop the
# number-format override back to Normal so no stray style lingers on the cell.
$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = "76442781"
$ws.Range("A3").Style = "Normal"

$ws.Range("B3").Value = 4080
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 4080
$ws.Range("G3").Value = "Cash"
$ws.Range("H3").Value = "2025-08-20T08:53:29"
